$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.355.54"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.094.44"
$ws.Range("E3").Value = "  +4.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.88"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9994"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5229"
$ws.Range("E7").Value = "  +1.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4344"
$ws.Range("E8").Value = "  +1.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08824"
$ws.Range("E9").Value = "  +1.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.88"
$ws.Range("E10").Value = "  +7.94%  "
$ws.Range("E11").Value = "  +2.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.51"
$ws.Range("E12").Value = "  -0.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.088.81"
$ws.Range("E13").Value = "  +3.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.729"
$ws.Range("E14").Value = "  +2.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.753"
$ws.Range("E15").Value = "  +4.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.50"
$ws.Range("E16").Value = "  +2.06%  "
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001129"
$ws.Range("E18").Value = "  +1.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06636"
$ws.Range("E19").Value = "  +1.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.95"
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9994"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.347"
$ws.Range("E22").Value = "  +2.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.393.11"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.37"
$ws.Range("E24").Value = "  +4.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.310"
$ws.Range("E25").Value = "  +2.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.333.52"
$ws.Range("E26").Value = "  +3.92%  "
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("E28").Value = "  +7.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "161.74"
$ws.Range("E29").Value = "  -0.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.18"
$ws.Range("E30").Value = "  +0.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.208"
$ws.Range("E31").Value = "  +6.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1072"
$ws.Range("E32").Value = "  +1.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.689"
$ws.Range("E33").Value = "  +24.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.175"
$ws.Range("E34").Value = "  +1.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.898"
$ws.Range("E35").Value = "  +1.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.998"
$ws.Range("E36").Value = "  +10.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02585"
$ws.Range("E37").Value = "  +1.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06701"
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.478"
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.67"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2260"
$ws.Range("E41").Value = "  +3.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6824"
$ws.Range("E42").Value = "  +2.71%  "
$ws.Range("E43").Value = "  +1.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9993"
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.09"
$ws.Range("E45").Value = "  +3.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6375"
$ws.Range("E46").Value = "  +3.45%  "
$ws.Range("E47").Value = "  +0.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.617"
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.193"
$ws.Range("E50").Value = "  +7.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.97"
$ws.Range("E51").Value = "  +1.42%  "
